# Bourbon_Chasers leaderboard - "latest. Updated wk 10"
# New weekly point totals came in; the standings (rows 2-10, sorted
# descending by Total Points) are refreshed/re-sorted accordingly.
# Row 11 (Brandon Greife, 182) is unaffected by this week's update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New standings after week 10, already in final (sorted) row order.
$ws.Range("A2").Value = "James Manley"
$ws.Range("B2").Value = 3291.5

$ws.Range("A3").Value = "Steven Carter"
$ws.Range("B3").Value = 3063.6

$ws.Range("A4").Value = "Matt Piunti"
$ws.Range("B4").Value = 2712

$ws.Range("A5").Value = "Jeremiah Gaddy"
$ws.Range("B5").Value = 2266.6999999999998

$ws.Range("A6").Value = "Todd Vinsant"
$ws.Range("B6").Value = 1936.4

$ws.Range("A7").Value = "Philip Milam"
$ws.Range("B7").Value = 1690.7

$ws.Range("A8").Value = "Senay Semere"
$ws.Range("B8").Value = 1655.3

$ws.Range("A9").Value = "Andrew Harrell"
$ws.Range("B9").Value = 1140.9000000000001

$ws.Range("A10").Value = "Josh Lance"
$ws.Range("B10").Value = 1101.7

# Row 11 (Brandon Greife, 182) is unchanged.

# Match the saved selection: the whole of row 11 selected (A11:XFD11),
# active cell A11.
$ws.Range("A11:XFD11").Select()

# Reflect the workbook window being moved/resized on save.
$excel.ActiveWindow.Left = -110
$excel.ActiveWindow.Top = -110
$excel.ActiveWindow.Width = 19420
$excel.ActiveWindow.Height = 10420
